$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F5").Value = 3076
    $ws.Range("F7").Value = 2423
    $ws.Range("F8").Value = 186
    $ws.Range("F16").Value = 297
}
